$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 452, pushing the old rows 452..533 down to 453..534.
$ws.Rows.Item(452).Insert()

# Populate the newly-inserted row 452 with the new weekly record.
$ws.Range("A452").Value = 4
$ws.Range("B452").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C452").Value = "Los Lagos"
$ws.Range("D452").Value = 45180
$ws.Range("E452").Value = 10
$ws.Range("F452").Value = 100112045
$ws.Range("G452").Value = "Zapallo"
$ws.Range("H452").Value = "Paine"
$ws.Range("I452").Value = "1a (guarda)"
$ws.Range("J452").Value = 500
$ws.Range("K452").Value = 800
$ws.Range("L452").Value = 800
$ws.Range("M452").Value = 800
$ws.Range("N452").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O452").Value = "Regi" + [char]0x00F3 + "n de O'Higgins"
$ws.Range("P452").Value = 800
$ws.Range("Q452").Value = 1
$ws.Range("R452").Value = "Hortaliza"
